{"js": "// Update the worksheet date and every \"a\u00d7b=c\" answer cell in the\n// multiplication-practice table to the new day's values.\n// Each entry is [oldText, newText] (old text is unique in the document).\nconst replacements = [\n  [\"2025-04-23 Wednesday\", \"2025-04-24 Thursday\"],\n  [\"61\u00d730=1830\", \"37\u00d740=1480\"],\n  [\"98\u00d760=5880\", \"87\u00d713=1131\"],\n  [\"32\u00d763=2016\", \"91\u00d793=8463\"],\n  [\"85\u00d796=8160\", \"27\u00d758=1566\"],\n  [\"88\u00d723=2024\", \"17\u00d779=1343\"],\n  [\"15\u00d733=495\", \"89\u00d775=6675\"],\n  [\"84\u00d774=6216\", \"94\u00d718=1692\"],\n  [\"20\u00d719=380\", \"97\u00d753=5141\"],\n  [\"18\u00d717=306\", \"44\u00d761=2684\"],\n  [\"76\u00d719=1444\", \"53\u00d749=2597\"],\n  [\"43\u00d766=2838\", \"64\u00d768=4352\"],\n  [\"23\u00d725=575\", \"31\u00d792=2852\"],\n  [\"64\u00d749=3136\", \"76\u00d776=5776\"],\n  [\"27\u00d793=2511\", \"85\u00d736=3060\"],\n  [\"59\u00d741=2419\", \"80\u00d731=2480\"],\n  [\"85\u00d787=7395\", \"29\u00d773=2117\"],\n  [\"40\u00d747=1880\", \"21\u00d792=1932\"],\n  [\"49\u00d735=1715\", \"55\u00d795=5225\"],\n  [\"31\u00d727=837\", \"63\u00d793=5859\"],\n  [\"42\u00d778=3276\", \"81\u00d754=4374\"],\n  [\"65\u00d737=2405\", \"11\u00d724=264\"],\n  [\"75\u00d779=5925\", \"48\u00d780=3840\"],\n  [\"81\u00d773=5913\", \"76\u00d785=6460\"],\n  [\"75\u00d739=2925\", \"20\u00d713=260\"],\n  [\"39\u00d724=936\", \"65\u00d713=845\"],\n];\n\nconst body = context.document.body;\n\n// Queue up a search for every old value, then sync once so all the\n// search results get populated together.\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nawait context.sync();\n\n// Replace each match with its corresponding new value, preserving the\n// run's existing formatting (insertText with Replace keeps the run's\n// rPr intact since it overwrites the text of the matched range).\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const results = searchResults[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every \"a\u00d7b=c\" answer cell in the\n# multiplication-practice table to the new day's values.\n# Each tuple is (oldText, newText) - old text is unique in the document.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-04-23 Wednesday', '2025-04-24 Thursday'),\n    @('61\u00d730=1830', '37\u00d740=1480'),\n    @('98\u00d760=5880', '87\u00d713=1131'),\n    @('32\u00d763=2016', '91\u00d793=8463'),\n    @('85\u00d796=8160', '27\u00d758=1566'),\n    @('88\u00d723=2024', '17\u00d779=1343'),\n    @('15\u00d733=495', '89\u00d775=6675'),\n    @('84\u00d774=6216', '94\u00d718=1692'),\n    @('20\u00d719=380', '97\u00d753=5141'),\n    @('18\u00d717=306', '44\u00d761=2684'),\n    @('76\u00d719=1444', '53\u00d749=2597'),\n    @('43\u00d766=2838', '64\u00d768=4352'),\n    @('23\u00d725=575', '31\u00d792=2852'),\n    @('64\u00d749=3136', '76\u00d776=5776'),\n    @('27\u00d793=2511', '85\u00d736=3060'),\n    @('59\u00d741=2419', '80\u00d731=2480'),\n    @('85\u00d787=7395', '29\u00d773=2117'),\n    @('40\u00d747=1880', '21\u00d792=1932'),\n    @('49\u00d735=1715', '55\u00d795=5225'),\n    @('31\u00d727=837', '63\u00d793=5859'),\n    @('42\u00d778=3276', '81\u00d754=4374'),\n    @('65\u00d737=2405', '11\u00d724=264'),\n    @('75\u00d779=5925', '48\u00d780=3840'),\n    @('81\u00d773=5913', '76\u00d785=6460'),\n    @('75\u00d739=2925', '20\u00d713=260'),\n    @('39\u00d724=936', '65\u00d713=845')\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
